$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds values in rows 3, 4, 6, 8 ---

# Row 3 (Corinthians - Cruzeiro)
$ws.Range("G3").Value = 1.38
$ws.Range("J3").Value = 1.91
$ws.Range("L3").Value = 8
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.9
$ws.Range("AJ3").Value = 23
$ws.Range("AO3").Value = 6.5

# Row 4 (Athletico-PR - Atletico GO)
$ws.Range("Q4").Value = 2.03
$ws.Range("R4").Value = 1.87

# Row 6 (Criciuma - Vitoria)
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62

# Row 8 (Eldense - Huesca)
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 5.5
$ws.Range("S8").Value = 1.73
$ws.Range("T8").Value = 2.08

# --- Append new row 9 with match data ---
$row9 = @("6uOnIaCm", "20/11/2024", "16:45", "WALES - CYMRU PREMIER", "Briton Ferry", "TNS", 40, 8.75, 1.04, 29, 3.5, 1.27, 1.02, 10, 1.08, 6.6, 1.26, 3.5, 1.18, 4.25, 2.95, 1.35, 200, 800, 300, 1000, 500, 800, 21, 32, 90, 500, 900, 11.5, 6.7, 17.5, 5.6, 14.5, 70, 40, 500, 250, "", 500, "", 4.25, 16, 175, 3.05, 3.75, 19, 6.5, 32, 300, "", "")

for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}
